$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7456
$ws1.Range("F4").Value = 279
$ws1.Range("F6").Value = 4025
$ws1.Range("F7").Value = 321
$ws1.Range("F8").Value = 567
$ws1.Range("F10").Value = 643
$ws1.Range("F11").Value = 126

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 6

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7456
$ws4.Range("F6").Value = 279
$ws4.Range("F8").Value = 4025
$ws4.Range("F9").Value = 321
$ws4.Range("F10").Value = 567
$ws4.Range("F12").Value = 643
$ws4.Range("F13").Value = 6
$ws4.Range("F14").Value = 126
